# Changed rounding for equity adjustment: np.ceil -> np.round
# Apply the resulting value changes to the "proposed_rhna_allocation" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("proposed_rhna_allocation")

$ws.Range("F2").Value = 1421
$ws.Range("C8").Value = 617
$ws.Range("F9").Value = 1317
$ws.Range("F11").Value = 6511
$ws.Range("F15").Value = 1251
$ws.Range("C20").Value = 744
$ws.Range("C23").Value = 198
$ws.Range("F23").Value = 344
$ws.Range("C31").Value = 326
$ws.Range("C41").Value = 151
$ws.Range("F41").Value = 262
$ws.Range("C43").Value = 20
$ws.Range("C48").Value = 634
$ws.Range("F52").Value = 103
$ws.Range("C53").Value = 213
$ws.Range("C55").Value = 12014
$ws.Range("C71").Value = 405
$ws.Range("F84").Value = 53
$ws.Range("C86").Value = 1597
$ws.Range("C88").Value = 8687
$ws.Range("F88").Value = 15088
$ws.Range("F92").Value = 828
$ws.Range("C95").Value = 447
$ws.Range("C106").Value = 701
$ws.Range("C110").Value = 222
